# Updates the cryptocurrency price/volume data to the latest scraped
# values (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.906.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.634.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.18'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.77%  '

# Row 9
$ws.Range("E9").Value = '  -3.40%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0611'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0878'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.866.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.04%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.635.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.90%  '

# Row 14
$ws.Range("E14").Value = '  -0.91%  '

# Row 15
$ws.Range("E15").Value = '  +0.44%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.95%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.911.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.30%  '

# Row 19
$ws.Range("E19").Value = '  -0.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.53%  '

# Row 21
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("E22").Value = '  -0.88%  '

# Row 23
$ws.Range("E23").Value = '  -3.45%  '

# Row 24
$ws.Range("E24").Value = '  -3.91%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.96%  '

# Row 28
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.111'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.83%  '

# Row 29
$ws.Range("E29").Value = '  +0.03%  '

# Row 30
$ws.Range("E30").Value = '  -1.29%  '

# Row 31
$ws.Range("E31").Value = '  -0.72%  '

# Row 32
$ws.Range("E32").Value = '  +0.56%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.397.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.94%  '

# Row 34
$ws.Range("E34").Value = '  -1.60%  '

# Row 35
$ws.Range("E35").Value = '  -0.40%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.79%  '

# Row 37
$ws.Range("E37").Value = '  +1.47%  '

# Row 38
$ws.Range("E38").Value = '  +0.13%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.558'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.870'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.58%  '

# Row 42
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.83%  '

# Row 44
$ws.Range("E44").Value = '  +2.62%  '

# Row 45
$ws.Range("E45").Value = '  -0.10%  '

# Row 46
$ws.Range("E46").Value = '  -1.78%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.775.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.29%  '

# Row 49
$ws.Range("E49").Value = '  -0.91%  '

# Row 50
$ws.Range("E50").Value = '  -0.24%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.97%  '
